$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# F-exam solutions posted -> per-student Final Exam raw scores (column N) filled in
$ws.Range("N12").Value = 128.0
$ws.Range("N13").Value = 120.0
$ws.Range("N14").Value = 112.0
$ws.Range("N15").Value = 104.0
$ws.Range("N16").Value = 95.0
$ws.Range("N17").Value = 85.0
$ws.Range("N18").Value = 74.0
$ws.Range("N19").Value = 64.0
$ws.Range("N20").Value = 54.0
$ws.Range("N21").Value = 44.0
$ws.Range("N22").Value = 34.0
$ws.Range("N23").Value = 0.0

# Final Exam score (F34) was graded and entered (1.0 = out of scale used elsewhere),
# and its raw total points (E34) is now populated with 141.
$ws.Range("E34").Value = 141.0
$ws.Range("F34").Value = 1.0
